$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Populate the new user rows (81-113): UserName, Password, Description, Locked, Email
$ws.Cells.Item(81, 1).Value = "KHPaddUser1"
$ws.Cells.Item(81, 2).Value = "Password1"
$ws.Cells.Item(81, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(81, 6).Value = "N"
$ws.Cells.Item(81, 7).Value = "KHPaddUser1@mailinator.com"

$ws.Cells.Item(82, 1).Value = "KHPaddUser2"
$ws.Cells.Item(82, 2).Value = "Password1"
$ws.Cells.Item(82, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(82, 6).Value = "N"
$ws.Cells.Item(82, 7).Value = "KHPaddUser2@mailinator.com"

$ws.Cells.Item(83, 1).Value = "KHPaddUser3"
$ws.Cells.Item(83, 2).Value = "Password1"
$ws.Cells.Item(83, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(83, 6).Value = "N"
$ws.Cells.Item(83, 7).Value = "KHPaddUser3@mailinator.com"

$ws.Cells.Item(84, 1).Value = "KHPaddUser4"
$ws.Cells.Item(84, 2).Value = "Password1"
$ws.Cells.Item(84, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(84, 6).Value = "N"
$ws.Cells.Item(84, 7).Value = "KHPaddUser4@mailinator.com"

$ws.Cells.Item(85, 1).Value = "KHPaddUser5"
$ws.Cells.Item(85, 2).Value = "Password1"
$ws.Cells.Item(85, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(85, 6).Value = "N"
$ws.Cells.Item(85, 7).Value = "KHPaddUser5@mailinator.com"

$ws.Cells.Item(86, 1).Value = "KHPaddUser6"
$ws.Cells.Item(86, 2).Value = "Password1"
$ws.Cells.Item(86, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(86, 6).Value = "N"
$ws.Cells.Item(86, 7).Value = "KHPaddUser6@mailinator.com"

$ws.Cells.Item(87, 1).Value = "SearchWhatsMarketUser1"
$ws.Cells.Item(87, 2).Value = "Password1"
$ws.Cells.Item(87, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(87, 6).Value = "N"
$ws.Cells.Item(87, 7).Value = "SearchWhatsMarketUser1@mailinator.com "

$ws.Cells.Item(88, 1).Value = "SearchWhatsMarketUser2"
$ws.Cells.Item(88, 2).Value = "Password1"
$ws.Cells.Item(88, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(88, 6).Value = "N"
$ws.Cells.Item(88, 7).Value = "SearchWhatsMarketUser2@mailinator.com "

$ws.Cells.Item(89, 1).Value = "SearchWhatsMarketUser3"
$ws.Cells.Item(89, 2).Value = "Password1"
$ws.Cells.Item(89, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(89, 6).Value = "N"
$ws.Cells.Item(89, 7).Value = "SearchWhatsMarketUser3@mailinator.com "

$ws.Cells.Item(90, 1).Value = "SearchWhatsMarketUser4"
$ws.Cells.Item(90, 2).Value = "Password1"
$ws.Cells.Item(90, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(90, 6).Value = "N"
$ws.Cells.Item(90, 7).Value = "SearchWhatsMarketUser4@mailinator.com "

$ws.Cells.Item(91, 1).Value = "SearchWhatsMarketUser5"
$ws.Cells.Item(91, 2).Value = "Password1"
$ws.Cells.Item(91, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(91, 6).Value = "N"
$ws.Cells.Item(91, 7).Value = "SearchWhatsMarketUser5@mailinator.com "

$ws.Cells.Item(92, 1).Value = "SearchWhatsMarketUser6"
$ws.Cells.Item(92, 2).Value = "Password1"
$ws.Cells.Item(92, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(92, 6).Value = "N"
$ws.Cells.Item(92, 7).Value = "SearchWhatsMarketUser6@mailinator.com "

$ws.Cells.Item(93, 1).Value = "SearchWhatsMarketUser7"
$ws.Cells.Item(93, 2).Value = "Password1"
$ws.Cells.Item(93, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(93, 6).Value = "N"
$ws.Cells.Item(93, 7).Value = "SearchWhatsMarketUser7@mailinator.com "

$ws.Cells.Item(94, 1).Value = "SearchWhatsMarketUser8"
$ws.Cells.Item(94, 2).Value = "Password1"
$ws.Cells.Item(94, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(94, 6).Value = "N"
$ws.Cells.Item(94, 7).Value = "SearchWhatsMarketUser8@mailinator.com "

$ws.Cells.Item(95, 1).Value = "SearchKnowHowUser1"
$ws.Cells.Item(95, 2).Value = "Password1"
$ws.Cells.Item(95, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(95, 6).Value = "N"
$ws.Cells.Item(95, 7).Value = "SearchKnowHowUser1@mailinator.com "

$ws.Cells.Item(96, 1).Value = "SearchKnowHowUser2"
$ws.Cells.Item(96, 2).Value = "Password1"
$ws.Cells.Item(96, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(96, 6).Value = "N"
$ws.Cells.Item(96, 7).Value = "SearchKnowHowUser2@mailinator.com "

$ws.Cells.Item(97, 1).Value = "SearchKnowHowUser3"
$ws.Cells.Item(97, 2).Value = "Password1"
$ws.Cells.Item(97, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(97, 6).Value = "N"
$ws.Cells.Item(97, 7).Value = "SearchKnowHowUser3@mailinator.com "

$ws.Cells.Item(98, 1).Value = "SearchKnowHowUser4"
$ws.Cells.Item(98, 2).Value = "Password1"
$ws.Cells.Item(98, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(98, 6).Value = "N"
$ws.Cells.Item(98, 7).Value = "SearchKnowHowUser4@mailinator.com "

$ws.Cells.Item(99, 1).Value = "SearchKnowHowUser5"
$ws.Cells.Item(99, 2).Value = "Password1"
$ws.Cells.Item(99, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(99, 6).Value = "N"
$ws.Cells.Item(99, 7).Value = "SearchKnowHowUser5@mailinator.com "

$ws.Cells.Item(100, 1).Value = "SearchKnowHowUser6"
$ws.Cells.Item(100, 2).Value = "Password1"
$ws.Cells.Item(100, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(100, 6).Value = "N"
$ws.Cells.Item(100, 7).Value = "SearchKnowHowUser6@mailinator.com "

$ws.Cells.Item(101, 1).Value = "SearchKnowHowUser7"
$ws.Cells.Item(101, 2).Value = "Password1"
$ws.Cells.Item(101, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(101, 6).Value = "N"
$ws.Cells.Item(101, 7).Value = "SearchKnowHowUser7@mailinator.com "

$ws.Cells.Item(102, 1).Value = "SearchKnowHowUser8"
$ws.Cells.Item(102, 2).Value = "Password1"
$ws.Cells.Item(102, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(102, 6).Value = "N"
$ws.Cells.Item(102, 7).Value = "SearchKnowHowUser8@mailinator.com "

$ws.Cells.Item(103, 1).Value = "AskUser1"
$ws.Cells.Item(103, 2).Value = "Password1"
$ws.Cells.Item(103, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(103, 6).Value = "N"
$ws.Cells.Item(103, 7).Value = "AskUser1@mailinator.com "

$ws.Cells.Item(104, 1).Value = "AskUser2"
$ws.Cells.Item(104, 2).Value = "Password1"
$ws.Cells.Item(104, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(104, 6).Value = "N"
$ws.Cells.Item(104, 7).Value = "AskUser2@mailinator.com "

$ws.Cells.Item(105, 1).Value = "AskUser3"
$ws.Cells.Item(105, 2).Value = "Password1"
$ws.Cells.Item(105, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(105, 6).Value = "N"
$ws.Cells.Item(105, 7).Value = "AskUser3@mailinator.com "

$ws.Cells.Item(106, 1).Value = "AskUser4"
$ws.Cells.Item(106, 2).Value = "Password1"
$ws.Cells.Item(106, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(106, 6).Value = "N"
$ws.Cells.Item(106, 7).Value = "AskUser4@mailinator.com "

$ws.Cells.Item(107, 1).Value = "AskUser5"
$ws.Cells.Item(107, 2).Value = "Password1"
$ws.Cells.Item(107, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(107, 6).Value = "N"
$ws.Cells.Item(107, 7).Value = "AskUser5@mailinator.com "

$ws.Cells.Item(108, 1).Value = "AskUser6"
$ws.Cells.Item(108, 2).Value = "Password1"
$ws.Cells.Item(108, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(108, 6).Value = "N"
$ws.Cells.Item(108, 7).Value = "AskUser6@mailinator.com "

$ws.Cells.Item(109, 1).Value = "AssetPageUser1"
$ws.Cells.Item(109, 2).Value = "Password1"
$ws.Cells.Item(109, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(109, 6).Value = "N"
$ws.Cells.Item(109, 7).Value = "AssetPageUser1@mailinator.com "

$ws.Cells.Item(110, 1).Value = "AssetPageUser2"
$ws.Cells.Item(110, 2).Value = "Password1"
$ws.Cells.Item(110, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(110, 6).Value = "N"
$ws.Cells.Item(110, 7).Value = "AssetPageUser2@mailinator.com "

$ws.Cells.Item(111, 1).Value = "AssetPageUser3"
$ws.Cells.Item(111, 2).Value = "Password1"
$ws.Cells.Item(111, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(111, 6).Value = "N"
$ws.Cells.Item(111, 7).Value = "AssetPageUser3@mailinator.com "

$ws.Cells.Item(112, 1).Value = "AssetPageUser4"
$ws.Cells.Item(112, 2).Value = "Password1"
$ws.Cells.Item(112, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(112, 6).Value = "N"
$ws.Cells.Item(112, 7).Value = "AssetPageUser4@mailinator.com "

$ws.Cells.Item(113, 1).Value = "AssetPageUser5"
$ws.Cells.Item(113, 2).Value = "Password1"
$ws.Cells.Item(113, 5).Value = "THIS IS IN USE 24/7 - DO NOT USE!"
$ws.Cells.Item(113, 6).Value = "N"
$ws.Cells.Item(113, 7).Value = "AssetPageUser5@mailinator.com "

# Apply the distinct Arial 10pt font used for the new UserName column cells
$ws.Range("A81:A113").Font.Name = "Arial"
$ws.Range("A81:A113").Font.Size = 10

# Add mailto hyperlinks on the Email column for the new rows (rows 81 and 112 are left without one, matching the source edit)
$ws.Hyperlinks.Add($ws.Cells.Item(82, 7), "mailto:KHPaddUser2@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(83, 7), "mailto:KHPaddUser3@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(84, 7), "mailto:KHPaddUser4@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(85, 7), "mailto:KHPaddUser5@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(86, 7), "mailto:KHPaddUser6@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(87, 7), "mailto:SearchWhatsMarketUser1@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(88, 7), "mailto:SearchWhatsMarketUser2@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(89, 7), "mailto:SearchWhatsMarketUser3@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(90, 7), "mailto:SearchWhatsMarketUser4@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(91, 7), "mailto:SearchWhatsMarketUser5@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(92, 7), "mailto:SearchWhatsMarketUser6@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(93, 7), "mailto:SearchWhatsMarketUser7@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(94, 7), "mailto:SearchWhatsMarketUser8@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(95, 7), "mailto:SearchKnowHowUser1@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(96, 7), "mailto:SearchKnowHowUser2@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(97, 7), "mailto:SearchKnowHowUser3@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(98, 7), "mailto:SearchKnowHowUser4@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(99, 7), "mailto:SearchKnowHowUser5@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(100, 7), "mailto:SearchKnowHowUser6@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(101, 7), "mailto:SearchKnowHowUser7@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(102, 7), "mailto:SearchKnowHowUser8@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(103, 7), "mailto:AskUser1@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(104, 7), "mailto:AskUser2@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(105, 7), "mailto:AskUser3@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(106, 7), "mailto:AskUser4@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(107, 7), "mailto:AskUser5@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(108, 7), "mailto:AskUser6@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(109, 7), "mailto:AssetPageUser1@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(110, 7), "mailto:AssetPageUser2@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(113, 7), "mailto:AssetPageUser5@mailinator.com")
$ws.Hyperlinks.Add($ws.Cells.Item(111, 7), "mailto:AssetPageUser3@mailinator.com", "", "", "AssetPageUser1@mailinator.com ")
$ws.Cells.Item(111, 7).Value = "AssetPageUser3@mailinator.com "
